$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Ativação date bump
# ----------------------------------------------------------------------
$d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2022", 2) | Out-Null

# Running paragraph index offset: original indices shift by +1 every time
# we split a paragraph into two (original content para + new italic para).
$offset = 0

# ----------------------------------------------------------------------
# 2) Objetivos (originally paragraph 6): replace PT text, append italic EN
# ----------------------------------------------------------------------
$idx = 6 + $offset
$pObjetivos = $d.Paragraphs($idx)
$pObjetivos.Range.Text = "Fornecer aos alunos conceitos fundamentais para compreensão da Química Inorgânica por meio da experimentação, desenvolvendo a capacidade de realizarem práticas no laboratório que estimulem o seu pensamento científico. Descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de caráter inorgânico com interesse industrial."
$pObjetivos.Range.InsertParagraphAfter()
$offset = $offset + 1

$idxEn = $idx + 1
$pObjetivosEn = $d.Paragraphs($idxEn)
$rObjetivosEn = $pObjetivosEn.Range
$enText1 = "Provide the students with fundamental concepts for understanding Inorganic Chemistry by means of experimentation, developing the ability to carry out practices in the laboratory which stimulate their scientific thinking. Describe and interpret the properties of the elements and their compounds, especially those of an inorganic feature with industrial interest."
$rObjetivosEn.Text = $enText1
$rIt1 = $d.Range($rObjetivosEn.Start, $rObjetivosEn.Start + $enText1.Length)
$rIt1.Italic = 1

# ----------------------------------------------------------------------
# 3) Programa resumido (originally paragraph 10): replace, append italic EN
# ----------------------------------------------------------------------
$idx = 10 + $offset
$pResumido = $d.Paragraphs($idx)
$pResumido.Range.Text = "Compostos de Coordenação. Materiais inorgânicos de interesse industrial. Purificação e Identificação de Compostos Inorgânicos. Síntese de sais e obtenção de Compostos de Alumínio."
$pResumido.Range.InsertParagraphAfter()
$offset = $offset + 1

$idxEn = $idx + 1
$pResumidoEn = $d.Paragraphs($idxEn)
$rResumidoEn = $pResumidoEn.Range
$enText2 = "Coordination Compounds. Inorganic materials of industrial interest. Purification and Identification of Inorganic Compounds. Synthesis: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds."
$rResumidoEn.Text = $enText2
$rIt2 = $d.Range($rResumidoEn.Start, $rResumidoEn.Start + $enText2.Length)
$rIt2.Italic = 1

# ----------------------------------------------------------------------
# 4) Programa (originally paragraph 12): replace, append italic EN
# ----------------------------------------------------------------------
$idx = 12 + $offset
$pPrograma = $d.Paragraphs($idx)
$pPrograma.Range.Text = "Compostos de Coordenação: Estrutura, ligações, reações e aplicações. Exemplos e aplicações de materiais inorgânicos de interesse industrial. Sínteses: Sal Simples, Sal Duplo e Sal Complexo. Preparação de Compostos de Alumínio."
$pPrograma.Range.InsertParagraphAfter()
$offset = $offset + 1

$idxEn = $idx + 1
$pProgramaEn = $d.Paragraphs($idxEn)
$rProgramaEn = $pProgramaEn.Range
$enText3 = "Coordination Compounds: Structure, bonds, reactions and applications. Examples and applications of industrial interest inorganic materials. Syntheses: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds."
$rProgramaEn.Text = $enText3
$rIt3 = $d.Range($rProgramaEn.Start, $rProgramaEn.Start + $enText3.Length)
$rIt3.Italic = 1

# ----------------------------------------------------------------------
# 5) Avaliação: Método / Critério / Norma de recuperação text swaps
# ----------------------------------------------------------------------
$d.Content.Find.Execute("P1 normal e P2 peso 2", $true, $false, $false, $false, $false, $true, 1, $false, "Serão oferecidas aulas expositivas e práticas.", 2) | Out-Null
$d.Content.Find.Execute("A média aritmética da prova P1 e P2 gera a nota final (NF)", $true, $false, $false, $false, $false, $true, 1, $false, "Serão aplicadas duas provas escritas. Trabalhos em sala de aula, seminários e relatórios, poderão, a critério do docente, ser considerados como parte da nota da prova escrita.", 2) | Out-Null
$d.Content.Find.Execute("Será fornecida uma aula na primeira semana seguida de uma prova escrita, na segunda semana a nota final recuperada será a média aritmética da NF com a da prova escrita.", $true, $false, $false, $false, $false, $true, 1, $false, "Será realizada uma prova escrita envolvendo o conteúdo do semestre todo.", 2) | Out-Null

# ----------------------------------------------------------------------
# 6) Bibliografia (originally paragraph 16): combine the 3 refs into one run
# ----------------------------------------------------------------------
# Use a whole-paragraph wildcard Find/Replace (rather than Range.Text=)
# so the merged run does not inherit the old first run's
# xml:space="preserve" (the old text began "1) Quagliano... 1973. " with
# a trailing space baked into its <w:t>).
$idx = 16 + $offset
$pBib = $d.Paragraphs($idx)
$rBib = $pBib.Range
$rBib.Find.Execute("1) Quagliano, J.V. e*Porto Alegre-RS, 2008.", $true, $false, $true, $false, $false, $true, 1, $false, "CHANG, Raymond. Química geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981.LEE, J. D., tradução Química Inorgânica não tão concisa da 5ª edição inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Química Inorgânica tradução da 4ª edição. Editora Bookman,Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Química - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3ª ed., 1973.", 2) | Out-Null
